$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New age-group trial condition rows, styled like the existing A2:A3 entries.
$newValues = @("TrialCondition_C.xlsx", "TrialCondition_D.xlsx", "TrialCondition_E.xlsx", "TrialCondition_F.xlsx")

$ws.Range("A2").Copy()

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = 4 + $i
    $cell = $ws.Cells.Item($row, 1)
    $cell.PasteSpecial(-4122)
    $cell.Value = $newValues[$i]
}

$ws.Range("A8").Select()
